$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new work-log entry as row 32, copying the formatting of the row
# directly above it (row 31) so the new cells share the same styles.
$ws.Range("A31:C31").Copy()
$ws.Range("A32:C32").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A32").Value = 45725
$ws.Range("B32").Value = 4
$ws.Range("C32").Value = "Struggled with integrating AI tool. It seems like there is problem with API key for AI tool. "

$ws.Range("C32").Select()
